# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table for the rows whose market data changed in this refresh.
# A value that looks like a plain decimal (e.g. '564.86') is written with a
# leading apostrophe so Excel keeps storing it as text, matching the
# original inlineStr/text cell type instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.964.70'
$ws.Range("E2").Value = '  +4.30%  '
$ws.Range("D3").Value = '2.465.60'
$ws.Range("E3").Value = '  +5.46%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''564.86'
$ws.Range("E5").Value = '  +2.91%  '
$ws.Range("D6").Value = '''142.72'
$ws.Range("E6").Value = '  +8.65%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''0.589'
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("D9").Value = '2.465.87'
$ws.Range("E9").Value = '  +5.49%  '
$ws.Range("E10").Value = '  +3.18%  '
$ws.Range("E11").Value = '  +1.51%  '
$ws.Range("E12").Value = '  +1.18%  '
$ws.Range("E13").Value = '  +4.17%  '
$ws.Range("D14").Value = '''26.31'
$ws.Range("E14").Value = '  +10.83%  '
$ws.Range("D15").Value = '2.906.41'
$ws.Range("E15").Value = '  +5.50%  '
$ws.Range("D16").Value = '62.880.43'
$ws.Range("E16").Value = '  +4.19%  '
$ws.Range("E17").Value = '  +4.39%  '
$ws.Range("D18").Value = '2.467.34'
$ws.Range("E18").Value = '  +6.15%  '
$ws.Range("E19").Value = '  +5.08%  '
$ws.Range("D20").Value = '''340.49'
$ws.Range("E20").Value = '  +7.99%  '
$ws.Range("D21").Value = '''4.26'
$ws.Range("E21").Value = '  +3.69%  '
$ws.Range("D22").Value = '''6.79'
$ws.Range("E22").Value = '  +2.77%  '
$ws.Range("D23").Value = '''0.998'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = '''65.46'
$ws.Range("E24").Value = '  +1.88%  '
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  +7.29%  '
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("E29").Value = '  +6.98%  '
$ws.Range("D30").Value = '''6.82'
$ws.Range("E30").Value = '  +11.19%  '
$ws.Range("D31").Value = '0.0₃0801'
$ws.Range("E31").Value = '  +9.19%  '
$ws.Range("E32").Value = '  +6.26%  '
$ws.Range("D33").Value = '''176.49'
$ws.Range("E33").Value = '  +3.04%  '
$ws.Range("E34").Value = '  +10.61%  '
$ws.Range("E35").Value = '  +3.09%  '
$ws.Range("D36").Value = '''18.85'
$ws.Range("E36").Value = '  +3.97%  '
$ws.Range("D37").Value = '''365.47'
$ws.Range("E37").Value = '  +12.60%  '
$ws.Range("D38").Value = '''4.39'
$ws.Range("E38").Value = '  +6.60%  '
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").Value = '''1.69'
$ws.Range("E41").Value = '  +10.06%  '
$ws.Range("D42").Value = '''40.47'
$ws.Range("E42").Value = '  +6.00%  '
$ws.Range("D43").Value = '''149.86'
$ws.Range("E43").Value = '  +8.70%  '
$ws.Range("D45").Value = '''20.52'
$ws.Range("E45").Value = '  +5.78%  '
$ws.Range("D46").Value = '''0.597'
$ws.Range("E46").Value = '  +5.09%  '
$ws.Range("D47").Value = '''0.0958'
$ws.Range("E47").Value = '  +0.73%  '
$ws.Range("D48").Value = '''0.0515'
$ws.Range("E48").Value = '  +3.03%  '
$ws.Range("D49").Value = '0.0₆0241'
$ws.Range("E49").Value = '  +9.70%  '
$ws.Range("E50").Value = '  +4.38%  '
$ws.Range("D51").Value = '''17.94'
$ws.Range("E51").Value = '  +4.77%  '
